# Loan RBI, Variable Instalments
# The "Repayment Schedule" sheet gains a new (blank) column at N, pushing
# the existing "Late" / (blank-heading) / "Outstanding" columns from
# N/O/P to O/P/Q respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N ("Late"), shifting N->O, O->P, P->Q.
$ws.Range("N1").EntireColumn.Insert()

# Author ended up with the selection on K19 after the edit.
$ws.Range("K19").Select()
